# The workbook gained one new daily price record (Mango, Vega Central
# Mapocho de Santiago) that belongs right above the existing row 219,
# pushing every subsequent record down by one row (219->220, ..., 328->329).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 219; Excel shifts rows 219:328 down to 220:329,
# carrying their values/styles with them.
$ws.Rows("219:219").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A219").Value2 = 9
$ws.Range("B219").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C219").Value2 = "Metropolitana"
$ws.Range("D219").Value2 = 44572
$ws.Range("E219").Value2 = 13
$ws.Range("F219").Value2 = "Fruta"
$ws.Range("G219").Value2 = 100108
$ws.Range("H219").Value2 = "Tropicales y subtropicales"
$ws.Range("I219").Value2 = 100108002
$ws.Range("J219").Value2 = "Mango"
$ws.Range("K219").Value2 = "Sin especificar"
$ws.Range("L219").Value2 = "Primera"
$ws.Range("M219").Value2 = 850
$ws.Range("N219").Value2 = 5500
$ws.Range("O219").Value2 = 6000
$ws.Range("P219").Value2 = 5741
$ws.Range("Q219").Value2 = "`$/bandeja 4 kilos"
$ws.Range("R219").Value2 = "Perú"
$ws.Range("S219").Value2 = 1435
$ws.Range("T219").Value2 = 4
